$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-10, leaving only row 1
$ws.Range("A2:A10").EntireRow.Delete()

# Update the remaining value in A1
$ws.Range("A1").Value = "Hi bro, bot test"

# Update selection to B6 as per the diff
$ws.Range("B6").Select()
